# Update countries & provincias Spain
# Applies the refreshed COVID-19 stats snapshot + two row re-ranks
# (Bahamas overtakes Somalia; Islas Malvinas overtakes Montserrat)
# plus the "last updated" timestamp on the Pais sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header: "last updated" timestamp (row 1, col A) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 23 de Septiembre de 2020 a las 03:09"

# --- Helper data: row => @(Pais, CasosTotales, NuevosCasos, CasosActivos, Recuperados, CasosCriticos, MuertesHoy, Muertes)
# Columns:            A      B             C            D            E            F            G          H

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 7097937
$ws.Cells.Item(4, 3).Value = 35696
$ws.Cells.Item(4, 4).Value = 4346110
$ws.Cells.Item(4, 5).Value = 2546356
$ws.Cells.Item(4, 7).Value = 969
$ws.Cells.Item(4, 8).Value = 205471

# Row 13: Argentina
$ws.Cells.Item(13, 2).Value = 652174
$ws.Cells.Item(13, 3).Value = 12027
$ws.Cells.Item(13, 4).Value = 517228
$ws.Cells.Item(13, 5).Value = 120994
$ws.Cells.Item(13, 7).Value = 470
$ws.Cells.Item(13, 8).Value = 13952

# Row 29: Canada
$ws.Cells.Item(29, 2).Value = 146663
$ws.Cells.Item(29, 3).Value = 1248
$ws.Cells.Item(29, 4).Value = 126904
$ws.Cells.Item(29, 5).Value = 10525

# Row 36: Panama
$ws.Cells.Item(36, 2).Value = 107284
$ws.Cells.Item(36, 3).Value = 474
$ws.Cells.Item(36, 4).Value = 83318
$ws.Cells.Item(36, 5).Value = 21681
$ws.Cells.Item(36, 7).Value = 13
$ws.Cells.Item(36, 8).Value = 2285

# Row 53: Venezuela
$ws.Cells.Item(53, 2).Value = 68453
$ws.Cells.Item(53, 3).Value = 1010
$ws.Cells.Item(53, 4).Value = 57774
$ws.Cells.Item(53, 5).Value = 10115
$ws.Cells.Item(53, 7).Value = 9
$ws.Cells.Item(53, 8).Value = 564

# Row 72: Paraguay
$ws.Cells.Item(72, 2).Value = 34828
$ws.Cells.Item(72, 3).Value = 568
$ws.Cells.Item(72, 4).Value = 19257
$ws.Cells.Item(72, 5).Value = 14866
$ws.Cells.Item(72, 7).Value = 29
$ws.Cells.Item(72, 8).Value = 705

# Row 125: Republica de Africa Central
$ws.Cells.Item(125, 2).Value = 4802
$ws.Cells.Item(125, 3).Value = 16
$ws.Cells.Item(125, 5).Value = 2910

# Rows 139/140: Bahamas overtakes Somalia in the ranking, so the two rows
# swap country names; Bahamas (row 139) gets the refreshed figures while
# Somalia (row 140) keeps its previous, unrevised figures.
$ws.Cells.Item(139, 1).Value = "Bahamas"
$ws.Cells.Item(139, 2).Value = 3467
$ws.Cells.Item(139, 3).Value = 49
$ws.Cells.Item(139, 4).Value = 1871
$ws.Cells.Item(139, 5).Value = 1519
$ws.Cells.Item(139, 7).Value = 2
$ws.Cells.Item(139, 8).Value = 77

$ws.Cells.Item(140, 1).Value = "Somalia"
$ws.Cells.Item(140, 2).Value = 3465
$ws.Cells.Item(140, 3).Value = 0
$ws.Cells.Item(140, 4).Value = 2877
$ws.Cells.Item(140, 5).Value = 490
$ws.Cells.Item(140, 7).Value = 0
$ws.Cells.Item(140, 8).Value = 98

# Row 186: Curazao
$ws.Cells.Item(186, 2).Value = 291
$ws.Cells.Item(186, 3).Value = 9
$ws.Cells.Item(186, 4).Value = 103
$ws.Cells.Item(186, 5).Value = 187

# Rows 214/215: Islas Malvinas overtakes Montserrat in the ranking, so the
# two rows fully swap their country name AND figures.
$ws.Cells.Item(214, 1).Value = "Islas Malvinas"
$ws.Cells.Item(214, 4).Value = 13
$ws.Cells.Item(214, 8).Value = 0

$ws.Cells.Item(215, 1).Value = "Montserrat"
$ws.Cells.Item(215, 4).Value = 12
$ws.Cells.Item(215, 8).Value = 1
